$wb = $excel.ActiveWorkbook

# Sheet "展览" - update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3462
$ws1.Range("F3").Value = 31
$ws1.Range("F5").Value = 1912
$ws1.Range("F6").Value = 137
$ws1.Range("F7").Value = 345

# Sheet "全部类型" - update "想去人数" (F column) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3462
$ws4.Range("F3").Value = 31
$ws4.Range("F5").Value = 1912
$ws4.Range("F6").Value = 137
$ws4.Range("F8").Value = 345
